$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. First paragraph (header line): change alignment from right to left
# ------------------------------------------------------------------
$d.Paragraphs(1).Range.ParagraphFormat.Alignment = 0

# ------------------------------------------------------------------
# 2. First table: shrink width from 5385 dxa (269.25 pt) to 5325 dxa (266.25 pt)
#    Need to touch tblPr width, the grid column and every cell in the
#    single-column table.
# ------------------------------------------------------------------
$t = $d.Tables(1)
$t.PreferredWidth = 266.25
$t.Columns(1).Width = 266.25
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cell = $t.Cell($i, 1)
    $cell.Width = 266.25
}

# ------------------------------------------------------------------
# 3. Swap the two TOC bookmarks wrapping the "Заявление о ..." heading
#    so that _Toc405368347 now comes first (id 0) and _Toc373232731
#    comes second (id 1). Bookmark ids are assigned by document order
#    at save time, so we delete the first one and re-add it at the
#    same spot - which places it after the other bookmark.
# ------------------------------------------------------------------
$bm1 = $d.Bookmarks("_Toc373232731")
$bmRange = $d.Range($bm1.Start, $bm1.End)
$bm1.Delete()
$d.Bookmarks.Add("_Toc373232731", $bmRange)

# ------------------------------------------------------------------
# 4. Merge the "202" + "1" runs into a single "2021" run (both runs
#    already share identical character formatting, so this is a pure
#    text/run merge). Re-assigning the very same text is treated as a
#    no-op by the engine and leaves the two runs intact, so the text
#    is first swapped to a placeholder and then back to "2021"; this
#    forces the range to be rewritten as a single run.
# ------------------------------------------------------------------
$yearRange = $d.Content
$yearRange.Find.Execute("2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$yearRange.Text = "PLACEHOLDER_YEAR"

$yearRange2 = $d.Content
$yearRange2.Find.Execute("PLACEHOLDER_YEAR", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)
$yearRange2.Text = "2021"
